$wb = $excel.ActiveWorkbook

# Each worksheet gets exactly one new row appended at the bottom, continuing
# the existing pattern of logged entries (next day's data, 2025-04-12).

$sheetsData = @(
    @{
        Name = "ROW50-FE-LIFTER"
        Row  = 70
        A    = 45759.24631425926
        B    = "0x01,0x90"
        C    = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D    = "0x01,0x52"
        E    = "0xe"
        F    = 400
        G    = [double]"5.68631262647114e+23"
        H    = 338
        I    = 14
    },
    @{
        Name = "ROW50-MID-LIFTER"
        Row  = 72
        A    = 45759.20769675926
        B    = "0x01,0x90 "
        C    = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D    = "0x01,0x56"
        E    = "0x19"
        F    = 400
        G    = "568631262647113771663628"
        H    = 342
        I    = 25
    },
    @{
        Name = "ROW11-FE-LIFTER"
        Row  = 70
        A    = 45759.27850303241
        B    = "0x01,0x90"
        C    = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D    = "0x01,0x52"
        E    = "0x14"
        F    = 400
        G    = [double]"5.68631262647114e+23"
        H    = 338
        I    = 20
    },
    @{
        Name = "ROW11-MID-LIFTER"
        Row  = 70
        A    = 45759.40496234954
        B    = "0x01,0x90"
        C    = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D    = "0x01,0x5a"
        E    = "0x19"
        F    = 400
        G    = [double]"5.68631262647114e+23"
        H    = 346
        I    = 25
    }
)

foreach ($entry in $sheetsData) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = $entry.F

    # Column G is normally numeric, but on the ROW50-MID-LIFTER sheet this
    # particular row stores a digit-string that is too large to round-trip
    # as a double, so it must be kept as text (matching the sheet's other
    # rows). Force text interpretation, then restore the default "Normal"
    # style so no stray number formatting is left on the cell.
    if ($entry.G -is [string]) {
        $gCell = $ws.Cells.Item($r, 7)
        $gCell.NumberFormat = "@"
        $gCell.Value = $entry.G
        $gCell.Style = "Normal"
    } else {
        $ws.Cells.Item($r, 7).Value = $entry.G
    }

    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
}
